$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Column A list: GOFAIPredictedGrowth/NeuralNetworkPredictedGrowth rows removed,
#    LastGDAXTrade moves up to A9.
$ws.Range("A9").Value = "LastGDAXTrade"
$ws.Range("A10").ClearContents() | Out-Null
$ws.Range("A11").ClearContents() | Out-Null

# 2. Remove the GOFAIPredictedGrowth and NeuralNetworkPredictedGrowth rows from the
#    normalisation table (old rows 12 & 13), shifting everything below up by two rows.
$ws.Rows("12:13").Delete() | Out-Null

# 3. Fix wording in the "Reasoning" column text (shared strings).
$ws.Range("K3").Value = "Top 100 currencies have less than or equal to 5 chars."
$ws.Range("K4").Value = "Top 100 currencies have less than or equal to 21 chars."

# 4. Update the active selection to match the saved view.
$ws.Range("K9").Select() | Out-Null
